$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (shared-string runs)
# ---------------------------------------------------------------------
# A8: "Volume 30   Number  24" -> "...25"
$ws.Range("A8").Value = "Volume 30   Number  25"
# C9: "Report Covering the Week  6/12/2023  Through  6/18/2023" -> new dates
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# ---------------------------------------------------------------------
# Helper: paste-formats-only constant
# xlPasteFormats = -4122
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
# C15: text "0" -> numeric 1  (style text -> number, copy format from F15 which stays numeric)
$ws.Range("C15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("F15").Value = 1

# G15: numeric 1 -> text "0" (style number -> text, copy format from D15 which stays text)
$ws.Range("G15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

# H15: numeric 100 -> text "***.*" (style number -> text, copy format from E15 which stays text)
$ws.Range("H15").Value = "***.*"
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -36.363636363636
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = -50

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = -55.555555555555
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 60
$ws.Range("H16").Value = -38.333333333333
$ws.Range("I16").Value = 237
$ws.Range("J16").Value = 280
$ws.Range("K16").Value = -15.357142857142
$ws.Range("L16").Value = 38.596491228070
$ws.Range("M16").Value = 264.615384615385
$ws.Range("N16").Value = -80.589680589680

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -10
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -41.860465116279
$ws.Range("I17").Value = 226
$ws.Range("J17").Value = 192
$ws.Range("K17").Value = 17.708333333333
$ws.Range("L17").Value = 5.116279069767
$ws.Range("M17").Value = 156.818181818182
$ws.Range("N17").Value = -27.564102564102

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 15
$ws.Range("E18").Value = -73.333333333333
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 65
$ws.Range("H18").Value = -61.538461538461
$ws.Range("I18").Value = 210
$ws.Range("J18").Value = 336
$ws.Range("K18").Value = -37.5
$ws.Range("L18").Value = 8.808290155440
$ws.Range("M18").Value = 28.834355828220
$ws.Range("N18").Value = -83.670295489891

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 50
$ws.Range("D19").Value = 57
$ws.Range("E19").Value = -12.280701754386
$ws.Range("F19").Value = 182
$ws.Range("G19").Value = 202
$ws.Range("H19").Value = -9.900990099009
$ws.Range("I19").Value = 1174
$ws.Range("J19").Value = 1008
$ws.Range("K19").Value = 16.468253968254
$ws.Range("L19").Value = 112.68115942029
$ws.Range("M19").Value = 6.727272727272
$ws.Range("N19").Value = -73.818019625334

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
# C20: text "0" -> numeric 1 (copy number format from D20, which stays numeric)
$ws.Range("C20").Value = 1
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -77.777777777777
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 39
$ws.Range("K20").Value = -30.769230769230
$ws.Range("L20").Value = 17.391304347826
$ws.Range("M20").Value = 107.692307692308
$ws.Range("N20").Value = -85.082872928176

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 73
$ws.Range("D21").Value = 101
$ws.Range("E21").Value = -27.722772277227
$ws.Range("F21").Value = 273
$ws.Range("G21").Value = 379
$ws.Range("H21").Value = -27.968337730870
$ws.Range("I21").Value = 1883
$ws.Range("J21").Value = 1871
$ws.Range("K21").Value = 0.641368252271
$ws.Range("L21").Value = 61.769759450171
$ws.Range("M21").Value = 31.586303284416
$ws.Range("N21").Value = -74.906716417910

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 14
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 111
$ws.Range("J22").Value = 92
$ws.Range("K22").Value = 20.652173913043
$ws.Range("L22").Value = 76.190476190476
$ws.Range("M22").Value = 76.190476190476
# N22 unchanged (text "***.*")

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 94
$ws.Range("D24").Value = 80
$ws.Range("E24").Value = 17.5
$ws.Range("F24").Value = 387
$ws.Range("G24").Value = 283
$ws.Range("H24").Value = 36.749116607773
$ws.Range("I24").Value = 1918
$ws.Range("J24").Value = 1451
$ws.Range("K24").Value = 32.184700206754
$ws.Range("L24").Value = 91.417165668662
$ws.Range("M24").Value = -15.357458075904
# N24 unchanged (text "***.*")

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 90
$ws.Range("G25").Value = 82
$ws.Range("H25").Value = 9.756097560975
$ws.Range("I25").Value = 490
$ws.Range("J25").Value = 427
$ws.Range("K25").Value = 14.754098360655
$ws.Range("L25").Value = 18.932038834951
$ws.Range("M25").Value = 74.377224199288
# N25 unchanged (text "***.*")

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
# C26: text "0" -> numeric 1 (copy number format from I26, which stays numeric)
$ws.Range("C26").Value = 1
$ws.Range("I26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("F26").Value = 2

# G26: numeric 2 -> text "0" (copy format from D26 which stays text)
$ws.Range("G26").Value = "'0"
$ws.Range("D26").Copy()
$ws.Range("G26").PasteSpecial(-4122)

# H26: numeric 50 -> text "***.*" (copy format from E26 which stays text)
$ws.Range("H26").Value = "***.*"
$ws.Range("E26").Copy()
$ws.Range("H26").PasteSpecial(-4122)

$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 13
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -27.777777777777
# M26, N26 unchanged (text "***.*")

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -62.5
$ws.Range("F27").Value = 21
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = -12.5
$ws.Range("I27").Value = 109
$ws.Range("J27").Value = 104
$ws.Range("K27").Value = 4.807692307692
$ws.Range("L27").Value = 62.686567164179
# M27, N27 unchanged (text "***.*")

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic. (only L28 changes)
# ---------------------------------------------------------------------
$ws.Range("L28").Value = -62.5

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc. (only L29 changes)
# ---------------------------------------------------------------------
$ws.Range("L29").Value = -66.666666666666

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------
# F30: numeric 1 -> text "0" (copy format from C30 which stays text)
$ws.Range("F30").Value = "'0"
$ws.Range("C30").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$ws.Range("H30").Value = -100
$ws.Range("L30").Value = -68.421052631578
